$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): force text format so values like "43.50" or "0.0420"
# are not coerced into numbers, then strip the temporary format so the
# cell ends up with no explicit style, matching a plain inline/shared text cell.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D11', 'D12', 'D14', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D34', 'D36', 'D37', 'D38', 'D40', 'D42', 'D45', 'D47', 'D48', 'D49')
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range('D2').Value = '93.400.70'
$ws.Range('D3').Value = '3.462.04'
$ws.Range('D5').Value = '235.09'
$ws.Range('D6').Value = '625.34'
$ws.Range('D11').Value = '3.460.24'
$ws.Range('D12').Value = '43.50'
$ws.Range('D14').Value = '6.24'
$ws.Range('D15').Value = '4.105.07'
$ws.Range('D16').Value = '93.262.36'
$ws.Range('D18').Value = '8.29'
$ws.Range('D19').Value = '3.457.51'
$ws.Range('D20').Value = '18.08'
$ws.Range('D21').Value = '11.81'
$ws.Range('D22').Value = '0.498'
$ws.Range('D23').Value = '3.40'
$ws.Range('D24').Value = '503.27'
$ws.Range('D25').Value = '6.80'
$ws.Range('D27').Value = '94.88'
$ws.Range('D28').Value = '12.20'
$ws.Range('D29').Value = '3.643.60'
$ws.Range('D30').Value = '2.84'
$ws.Range('D31').Value = '11.34'
$ws.Range('D34').Value = '1.02'
$ws.Range('D36').Value = '29.56'
$ws.Range('D37').Value = '0.556'
$ws.Range('D38').Value = '572.72'
$ws.Range('D40').Value = '7.53'
$ws.Range('D42').Value = '0.917'
$ws.Range('D45').Value = '0.0420'
$ws.Range('D47').Value = '5.50'
$ws.Range('D48').Value = '3.56'
$ws.Range('D49').Value = '53.20'

foreach ($c in $priceCells) { $ws.Range($c).ClearFormats() }

# Volume(1h) column (E): values always contain "%" and padding spaces,
# so Excel keeps them as text automatically.
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('E3').Value = '  +4.16%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +3.28%  '
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  +7.59%  '
$ws.Range('E8').Value = '  +3.58%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('E10').Value = '  +11.43%  '
$ws.Range('E11').Value = '  +4.13%  '
$ws.Range('E12').Value = '  +9.20%  '
$ws.Range('E13').Value = '  +5.34%  '
$ws.Range('E14').Value = '  +7.19%  '
$ws.Range('E15').Value = '  +3.99%  '
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('E17').Value = '  +3.41%  '
$ws.Range('E18').Value = '  +5.59%  '
$ws.Range('E19').Value = '  +4.22%  '
$ws.Range('E20').Value = '  +9.25%  '
$ws.Range('E21').Value = '  +10.27%  '
$ws.Range('E22').Value = '  +13.78%  '
$ws.Range('E23').Value = '  +10.83%  '
$ws.Range('E24').Value = '  +3.87%  '
$ws.Range('E25').Value = '  +11.06%  '
$ws.Range('E26').Value = '  +1.99%  '
$ws.Range('E27').Value = '  +7.24%  '
$ws.Range('E28').Value = '  +8.43%  '
$ws.Range('E29').Value = '  +4.62%  '
$ws.Range('E30').Value = '  +10.23%  '
$ws.Range('E31').Value = '  +2.58%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  +7.25%  '
$ws.Range('E34').Value = '  +2.74%  '
$ws.Range('E35').Value = '  +7.37%  '
$ws.Range('E36').Value = '  +5.33%  '
$ws.Range('E37').Value = '  +8.12%  '
$ws.Range('E38').Value = '  +10.66%  '
$ws.Range('E39').Value = '  +6.63%  '
$ws.Range('E40').Value = '  +3.46%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +7.02%  '
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('E45').Value = '  +9.02%  '
$ws.Range('E46').Value = '  +3.47%  '
$ws.Range('E47').Value = '  +3.69%  '
$ws.Range('E48').Value = '  -0.84%  '
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('E51').Value = '  +4.80%  '

